# Update countries & provincias Spain
# - Re-sorted several countries (by total cases) which changes which
#   country name lands on a given data row.
# - Refreshed the numeric Covid-19 figures (Casos totales, Nuevos casos,
#   Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
#   for the affected rows.
# - Bumped the "Datos actualizados..." timestamp banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Country name (column A) re-ordering caused by the new sort order.
#    The B:H figures for these rows are rewritten further below.
# ---------------------------------------------------------------------

# Peru moves above Republica Dominicana / Mexico
$ws.Range("A45").Value = "Peru"
$ws.Range("A46").Value = "Republica Dominicana"
$ws.Range("A47").Value = "Mexico"

# Kazajistan moves above Principado de Andorra / Eslovaquia
$ws.Range("A77").Value = "Kazajistan"
$ws.Range("A78").Value = "Principado de Andorra"
$ws.Range("A79").Value = "Eslovaquia"

# Togo moves above Polinesia Francesa / Mali
$ws.Range("A143").Value = "Togo"
$ws.Range("A144").Value = "Polinesia Francesa"
$ws.Range("A145").Value = "Mali"

# Nueva Caledonia moves above Gabon
$ws.Range("A158").Value = "Nueva Caledonia"
$ws.Range("A159").Value = "Gabon"

# Benin moves above Santa Lucia
$ws.Range("A167").Value = "Benin"
$ws.Range("A168").Value = "Santa Lucia"

# ---------------------------------------------------------------------
# 2) Updated figures (columns B:H) for the affected rows.
# ---------------------------------------------------------------------

# Estados Unidos
$ws.Range("B4").Value = 235281
$ws.Range("C4").Value = 20278
$ws.Range("D4").Value = 10324
$ws.Range("E4").Value = 219354
$ws.Range("G4").Value = 501
$ws.Range("H4").Value = 5603

# Alemania
$ws.Range("B7").Value = 83875
$ws.Range("C7").Value = 5894
$ws.Range("E7").Value = 61409
$ws.Range("G7").Value = 135
$ws.Range("H7").Value = 1066

# Francia
$ws.Range("B9").Value = 59105
$ws.Range("C9").Value = 2116
$ws.Range("D9").Value = 12428
$ws.Range("E9").Value = 42174
$ws.Range("G9").Value = 471
$ws.Range("H9").Value = 4503

# Australia
$ws.Range("B23").Value = 5139
$ws.Range("C23").Value = 91
$ws.Range("E23").Value = 4769

# Finlandia
$ws.Range("F43").Value = 65

# Sudafrica
$ws.Range("B44").Value = 1462
$ws.Range("C44").Value = 82
$ws.Range("E44").Value = 1407

# Peru (now row 45)
$ws.Range("B45").Value = 1414
$ws.Range("C45").Value = 91
$ws.Range("D45").Value = 394
$ws.Range("E45").Value = 973
$ws.Range("F45").Value = 49
$ws.Range("G45").Value = 9
$ws.Range("H45").Value = 47

# Republica Dominicana (now row 46)
$ws.Range("B46").Value = 1380
$ws.Range("C46").Value = 96
$ws.Range("D46").Value = 16
$ws.Range("E46").Value = 1304
$ws.Range("F46").Value = 147
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 60

# Mexico (now row 47)
$ws.Range("B47").Value = 1378
$ws.Range("C47").Value = 163
$ws.Range("D47").Value = 35
$ws.Range("E47").Value = 1306
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 37

# Barein
$ws.Range("B69").Value = 643
$ws.Range("C69").Value = 74
$ws.Range("D69").Value = 381
$ws.Range("E69").Value = 258

# Kazajistan (now row 77)
$ws.Range("B77").Value = 435
$ws.Range("C77").Value = 55
$ws.Range("D77").Value = 27
$ws.Range("E77").Value = 405
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 3

# Principado de Andorra (now row 78)
$ws.Range("B78").Value = 428
$ws.Range("C78").Value = 38
$ws.Range("D78").Value = 10
$ws.Range("E78").Value = 403
$ws.Range("F78").Value = 12
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 15

# Eslovaquia (now row 79)
$ws.Range("B79").Value = 426
$ws.Range("C79").Value = 26
$ws.Range("D79").Value = 5
$ws.Range("E79").Value = 420
$ws.Range("F79").Value = 3
$ws.Range("H79").Value = 1

# Togo (now row 143)
$ws.Range("B143").Value = 39
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 17
$ws.Range("E143").Value = 20
$ws.Range("F143").Value = 0
$ws.Range("H143").Value = 2

# Polinesia Francesa (now row 144)
$ws.Range("B144").Value = 37
$ws.Range("C144").Value = 0
$ws.Range("E144").Value = 37
$ws.Range("F144").Value = 1
$ws.Range("H144").Value = 0

# Mali (now row 145)
$ws.Range("C145").Value = 5
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 33
$ws.Range("H145").Value = 3

# Nueva Caledonia (now row 158)
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 1
$ws.Range("H158").Value = 0

# Gabon (now row 159)
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 0
$ws.Range("H159").Value = 1

# ---------------------------------------------------------------------
# 3) Timestamp banner update.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 19:50"
